# Refactor the "Shortest Path" / "Graph" / new "Linked List" sections:
# - apply the grey "background1 shaded A6" colour used elsewhere in the
#   document to the existing "Shortest Path" and "Graph" blocks
# - split "Graph:" into three runs and rename to "Graph (Undirected
#   Weighted):"
# - add a new "Linked List (Sorted Linked List):" block
# - move the <w:lastRenderedPageBreak/> marker from the "Cac con duong..."
#   bullet to the "Others" heading

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$grey = '<w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>'

function Get-ParaIndexByText($doc, $text) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -eq $text) {
            return $idx
        }
    }
    return -1
}

function Get-ParaByIndex($doc, $index) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($idx -eq $index) {
            return $p
        }
    }
    return $null
}

function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Shortest Path" heading (Heading2) -> add grey colour
# ---------------------------------------------------------------------
$t1 = "Shortest Path`r"
$p1 = Get-ParaByText $d $t1
$xml1 = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>Shortest Path</w:t></w:r></w:p>"
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "Dung de di chuyen theo kieu bam chuot vao man hinh " -> add grey colour
# ---------------------------------------------------------------------
$t2 = "Dùng để di chuyển theo kiểu bấm chuột vào màn hình `r"
$p2 = Get-ParaByText $d $t2
$xml2 = "<w:p $wNs><w:pPr><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t xml:space=`"preserve`">Dùng để di chuyển theo kiểu bấm chuột vào màn hình </w:t></w:r></w:p>"
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "(Tham khao: Xem lai bai tap Graph - Shortest Path)" -> add grey colour
# ---------------------------------------------------------------------
$t3 = "(Tham khảo: Xem lại bài tập Graph – Shortest Path)`r"
$p3 = Get-ParaByText $d $t3
$xml3 = "<w:p $wNs><w:pPr><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>(Tham khảo: Xem lại bài tập Graph – Shortest Path)</w:t></w:r></w:p>"
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) "Graph:" heading (Heading2) -> "Graph (Undirected Weighted):" (3 runs) + grey colour
# ---------------------------------------------------------------------
$t4 = "Graph:`r"
$p4 = Get-ParaByText $d $t4
$xml4 = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>Graph</w:t></w:r><w:r><w:rPr>$grey</w:rPr><w:t xml:space=`"preserve`"> (Undirected Weighted)</w:t></w:r><w:r><w:rPr>$grey</w:rPr><w:t>:</w:t></w:r></w:p>"
$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------
# 5) "Lam cung voi bai Shortest Path" (right after Graph) -> add grey colour
# ---------------------------------------------------------------------
$t5 = "Làm cùng với bài Shortest Path`r"
$p5 = Get-ParaByText $d $t5
$xml5 = "<w:p $wNs><w:pPr><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>Làm cùng với bài Shortest Path</w:t></w:r></w:p>"
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------
# 6) Insert two new paragraphs (Linked List heading + body) right before
#    the "File I/O" heading, i.e. right after the blank paragraph that
#    follows the Graph section.
#
#    InsertXML on a collapsed range merges the first/last paragraph of
#    the inserted fragment into the surrounding paragraphs (only runs
#    survive, pPr is dropped). To get two brand-new, fully-formatted
#    paragraphs we sandwich them between two throw-away empty
#    paragraphs and then delete the surplus empty paragraphs that
#    result, leaving exactly one blank line (matching the pre-existing
#    blank line) between the two sections.
# ---------------------------------------------------------------------
$tFileIO = "File I/O`r"
$idxFileIO = Get-ParaIndexByText $d $tFileIO
$idxBlank = $idxFileIO - 1
$pBlank = Get-ParaByIndex $d $idxBlank
$insPos = $pBlank.Range.Start
$insRng = $d.Range($insPos, $insPos)

$xml6 = "<w:p $wNs></w:p>" + `
        "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>Linked List (Sorted Linked List):</w:t></w:r></w:p>" + `
        "<w:p $wNs><w:pPr><w:rPr>$grey</w:rPr></w:pPr><w:r><w:rPr>$grey</w:rPr><w:t>Làm cùng với bài Shortest Path</w:t></w:r></w:p>" + `
        "<w:p $wNs></w:p>"
$insRng.InsertXML($xml6)

$tNewHeading = "Linked List (Sorted Linked List):`r"
$idxNewHeading = Get-ParaIndexByText $d $tNewHeading
$idxSurplusBefore = $idxNewHeading - 1
$pSurplusBefore = Get-ParaByIndex $d $idxSurplusBefore
$pSurplusBefore.Range.Delete()

$idxNewHeading2 = Get-ParaIndexByText $d $tNewHeading
$idxSurplusAfter = $idxNewHeading2 + 2
$pSurplusAfter = Get-ParaByIndex $d $idxSurplusAfter
$pSurplusAfter.Range.Delete()

# ---------------------------------------------------------------------
# 7) Move <w:lastRenderedPageBreak/>: remove from the "Cac con duong va o
#    dat chua mua..." bullet and add to the "Others" heading run.
# ---------------------------------------------------------------------
$t7 = "Các con đường và ô đất chưa mua sẽ chỉ là background (ô đất thì có collider)`r"
$p7 = Get-ParaByText $d $t7
$xml7 = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Các con đường và ô đất chưa mua sẽ chỉ là background (ô đất thì có collider)</w:t></w:r></w:p>"
$p7.Range.InsertXML($xml7)

$t8 = "Others`r"
$p8 = Get-ParaByText $d $t8
$xml8 = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Others</w:t></w:r></w:p>"
$p8.Range.InsertXML($xml8)

Write-Host "Edits applied successfully"
